$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2, 3, 4, 5, 7, 9, 10, 11, 12, 14, 15)

$data = @(
    @(0.963223312315705, 0.03974516942150075, 0.5059138127533629, 0.158705431593404, 0.002534080523199676, 1.415813774585239, 0.05490357979274352, 0.6249760980557255, 0.4624399593449908, 2.443596874654283, 6.118727321080826),
    @(0.9270537950513074, 0.03712838881840241, 0.5033804456136863, 0.1587916105771026, 0.002536787998756027, 1.423793513952614, 0.05511580796576698, 0.5900869330664023, 0.4582616829314361, 2.465075704722731, 6.145136798158575),
    @(0.9052260422830614, 0.03550481443962639, 0.5020343173756601, 0.1588988353522254, 0.002538540710204663, 1.4292311384879, 0.0552545306654304, 0.5688925569086791, 0.4558867513708122, 2.478934766275295, 6.163844719623853),
    @(0.8964274500786189, 0.03483896425977662, 0.5015385409784585, 0.1589562204919446, 0.002539277730847644, 1.431582428596176, 0.0553131819292485, 0.5603135005771946, 0.4549670069863794, 2.484751228659984, 6.172095643176903),
    @(0.894972293730433, 0.03472814510985955, 0.5014594095241875, 0.1589665770849891, 0.0025394014900611, 1.431981040778613, 0.0553230491533272, 0.5588924638488493, 0.4548171910745822, 2.485727243666556, 6.173503606288108),
    @(0.90510699014024, 0.03549585166906155, 0.5020274172721457, 0.1588995537928302, 0.002538550557600057, 1.429262300362446, 0.05525531306403497, 0.5687766218795502, 0.4558741525923722, 2.479012525648761, 6.163953453836911),
    @(0.9506735136136513, 0.03884640773859616, 0.5049969317174998, 0.1587238936792268, 0.002534995357389703, 1.418453592257833, 0.05497501276402073, 0.6128993787303898, 0.4609598262125445, 2.450863553226245, 6.127316342165301),
    @(1.043021858287176, 0.04528293060010924, 0.5124762930391142, 0.1588090054569555, 0.002528737127565443, 1.401521797431293, 0.05449188141749595, 0.7012107418910603, 0.472439308063187, 2.400983825674579, 6.075225345968988),
    @(1.11266680435125, 0.0499304546242314, 0.5189753306755875, 0.1591317789329807, 0.002524569886878254, 1.39167526861106, 0.05417717936496924, 0.7671624631603891, 0.4817858088972713, 2.367574525641791, 6.048971776663251),
    @(1.144734808905383, 0.05202713642208323, 0.5221487797031159, 0.159334790203733, 0.002522766706599976, 1.387757670768636, 0.05404268890287867, 0.797394012635749, 0.486234765032961, 2.353077432871054, 6.039632948654571),
    @(1.156933032974564, 0.05281857482189878, 0.5233815691999268, 0.1594197130513315, 0.00252209712454821, 1.386354836660473, 0.05399300247982186, 0.8088744779879278, 0.4879476982384716, 2.34768850386706, 6.036470581287347),
    @(1.154303503564876, 0.05264823710979272, 0.5231146860324429, 0.1594010660138814, 0.002522240742944915, 1.38665337542264, 0.05400364816278902, 0.8064005215563839, 0.4875775350831333, 2.348844622649054, 6.037135024554033),
    @(1.145737270909564, 0.05209229930574111, 0.5222495800378368, 0.1593416157674561, 0.002522711354838789, 1.387640642410837, 0.05403857630048758, 0.798337870896404, 0.4863751245376733, 2.352632062094587, 6.039365284756713),
    @(1.140497315635912, 0.05175144154603117, 0.5217237203776079, 0.1593062477763105, 0.002523001340789062, 1.388255875648944, 0.05406013243369046, 0.7934034703227724, 0.4856422828399616, 2.354965106227807, 6.040780082830736),
    @(1.110578790238691, 0.0497930782811693, 0.5187722937381096, 0.1591196394402772, 0.002524689585172242, 1.391942586918809, 0.05418614269968014, 0.7651913320649442, 0.4814990146384588, 2.368536050527004, 6.049634455650789),
    @(1.092323107538789, 0.04858719440973402, 0.5170171753582764, 0.1590195277172626, 0.002525748921082607, 1.394348051129306, 0.05426566312961789, 0.7479425184558011, 0.4790076546560442, 2.377040923632437, 6.055732990394176),
    @(1.081859315426698, 0.04789195548178782, 0.5160281110273246, 0.1589672359224288, 0.002526366935088098, 1.395784479685616, 0.0543122174138011, 0.7380431215090937, 0.4775932590402618, 2.381998679220558, 6.059485836642551),
    @(1.078322729151608, 0.04765627714688492, 0.515696744299035, 0.1589504403187227, 0.002526577682142084, 1.396279912728104, 0.05432812021888722, 0.7346950952233726, 0.4771175626077735, 2.38368862160992, 6.060798599991216),
    @(1.094262695047604, 0.04871573329911882, 0.5172018968866752, 0.1590296375347968, 0.002525635251769483, 1.394086514233031, 0.05425711359892826, 0.7497764463315661, 0.4792709434346563, 2.376128736898188, 6.055058423897748),
    @(1.148251900823482, 0.05225566046829044, 0.5225028403699952, 0.1593588596083038, 0.002522572765750508, 1.387348469358265, 0.05402828338206511, 0.8007051914464114, 0.4867275372043878, 2.351516863416055, 6.038700055119989),
    @(1.183855835735358, 0.05455445760954092, 0.5261483707697607, 0.1596209139940896, 0.002520648415698225, 1.383414964035097, 0.05388596805579571, 0.8341788362755267, 0.49176521489251, 2.336019158025985, 6.030188940045861),
    @(1.164824421535457, 0.05332890023819914, 0.5241861593443957, 0.1594767707807598, 0.002521668437290188, 1.38547135559871, 0.05396126355280639, 0.8162962510637044, 0.4890615202215258, 2.344236808486844, 6.034532146501846),
    @(1.09338570891137, 0.04865762697265552, 0.5171183220503366, 0.1590250504880331, 0.002525686613715956, 1.394204588421367, 0.05426097623433002, 0.748947273992087, 0.4791518547474283, 2.376540924106113, 6.055362626861637),
    @(1.01772167936096, 0.04355599086622419, 0.5102762128848326, 0.1587401609321937, 0.002530354204184289, 1.405646446858938, 0.05461548999927723, 0.6771309496192259, 0.4691731312010177, 2.413908677023045, 6.087205112400824)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $data[$i][$j]
    }
}
